$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format to avoid numeric auto-conversion / precision loss
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.321.21"
$ws.Range("E2").Value = "  -2.27%  "

$ws.Range("D3").Value = "3.477.43"
$ws.Range("E3").Value = "  -3.54%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "583.27"
$ws.Range("E5").Value = "  -3.62%  "

$ws.Range("D6").Value = "131.42"
$ws.Range("E6").Value = "  -4.30%  "

$ws.Range("D7").Value = "3.478.89"
$ws.Range("E7").Value = "  -3.50%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -2.06%  "

$ws.Range("E10").Value = "  -1.69%  "

$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("E12").Value = "  -1.76%  "

$ws.Range("D13").Value = "4.077.69"
$ws.Range("E13").Value = "  -3.35%  "

$ws.Range("D14").Value = "27.61"
$ws.Range("E14").Value = "  -1.77%  "

$ws.Range("E15").Value = "  -5.07%  "

$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("D17").Value = "3.485.50"
$ws.Range("E17").Value = "  -3.31%  "

$ws.Range("D18").Value = "64.395.26"
$ws.Range("E18").Value = "  -2.32%  "

$ws.Range("D19").Value = "9.89"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("E20").Value = "  -2.46%  "

$ws.Range("D21").Value = "5.65"
$ws.Range("E21").Value = "  -4.74%  "

$ws.Range("D22").Value = "392.11"
$ws.Range("E22").Value = "  -1.57%  "

$ws.Range("D23").Value = "0.574"
$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("D24").Value = "3.622.01"
$ws.Range("E24").Value = "  -3.46%  "

$ws.Range("D25").Value = "73.01"
$ws.Range("E25").Value = "  -2.11%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  -9.91%  "

$ws.Range("E28").Value = "  -7.63%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").Value = "7.38"
$ws.Range("E30").Value = "  -9.98%  "

$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  -6.75%  "

$ws.Range("D32").Value = "8.12"
$ws.Range("E32").Value = "  -6.32%  "

$ws.Range("D33").Value = "3.479.22"
$ws.Range("E33").Value = "  -3.57%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").Value = "23.90"
$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("D36").Value = "0.145"
$ws.Range("E36").Value = "  -2.16%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "5.19"
$ws.Range("E37").Value = "  -3.72%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "170.84"
$ws.Range("E38").Value = "  -0.66%  "

$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "6.95"
$ws.Range("E39").Value = "  -2.29%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "1.57"
$ws.Range("E40").Value = "  -2.89%  "

$ws.Range("D41").Value = "0.0803"
$ws.Range("E41").Value = "  -4.01%  "

$ws.Range("E42").Value = "  -3.69%  "

$ws.Range("D43").Value = "25.84"
$ws.Range("E43").Value = "  -1.05%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").Value = "41.79"
$ws.Range("E45").Value = "  -3.74%  "

$ws.Range("E46").Value = "  -4.40%  "

$ws.Range("E47").Value = "  -4.72%  "

$ws.Range("E48").Value = "  -4.63%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "6.86"
$ws.Range("E49").Value = "  -3.40%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.428.30"
$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "0.888"
$ws.Range("E51").Value = "  -0.17%  "
